$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-02-05 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-02-06 Monday", 2) | Out-Null

# Update the math-expression table cells directly by position (avoids Find/Replace
# substring collisions between already-updated and not-yet-updated cell text)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "61+5="
$tbl.Cell(1, 2).Range.Text = "43+44="
$tbl.Cell(1, 3).Range.Text = "78-13="
$tbl.Cell(1, 4).Range.Text = "34+62="
$tbl.Cell(1, 5).Range.Text = "50+47="
$tbl.Cell(2, 1).Range.Text = "92-39="
$tbl.Cell(2, 2).Range.Text = "20-2="
$tbl.Cell(2, 3).Range.Text = "69+13="
$tbl.Cell(2, 4).Range.Text = "28+60="
$tbl.Cell(2, 5).Range.Text = "51+30="
$tbl.Cell(3, 1).Range.Text = "38-8="
$tbl.Cell(3, 2).Range.Text = "16+30="
$tbl.Cell(3, 3).Range.Text = "78-48="
$tbl.Cell(3, 4).Range.Text = "7-5="
$tbl.Cell(3, 5).Range.Text = "76-43="
$tbl.Cell(4, 1).Range.Text = "43+11="
$tbl.Cell(4, 2).Range.Text = "11+63="
$tbl.Cell(4, 3).Range.Text = "42-37="
$tbl.Cell(4, 4).Range.Text = "66+13="
$tbl.Cell(4, 5).Range.Text = "59-30="
$tbl.Cell(5, 1).Range.Text = "34-3="
$tbl.Cell(5, 2).Range.Text = "98-36="
$tbl.Cell(5, 3).Range.Text = "11-2="
$tbl.Cell(5, 4).Range.Text = "88-82="
$tbl.Cell(5, 5).Range.Text = "40+56="
$tbl.Cell(6, 1).Range.Text = "39-14="
$tbl.Cell(6, 2).Range.Text = "96-88="
$tbl.Cell(6, 3).Range.Text = "55+31="
$tbl.Cell(6, 4).Range.Text = "40+5="
$tbl.Cell(6, 5).Range.Text = "27+48="
$tbl.Cell(7, 1).Range.Text = "27+21="
$tbl.Cell(7, 2).Range.Text = "73-33="
$tbl.Cell(7, 3).Range.Text = "62-39="
$tbl.Cell(7, 4).Range.Text = "34+11="
$tbl.Cell(7, 5).Range.Text = "99+0="
$tbl.Cell(8, 1).Range.Text = "85-81="
$tbl.Cell(8, 2).Range.Text = "93-56="
$tbl.Cell(8, 3).Range.Text = "81-27="
$tbl.Cell(8, 4).Range.Text = "46+2="
$tbl.Cell(8, 5).Range.Text = "42+39="
$tbl.Cell(9, 1).Range.Text = "0+56="
$tbl.Cell(9, 2).Range.Text = "18+46="
$tbl.Cell(9, 3).Range.Text = "28+29="
$tbl.Cell(9, 4).Range.Text = "5+24="
$tbl.Cell(9, 5).Range.Text = "64-13="
$tbl.Cell(10, 1).Range.Text = "97-50="
$tbl.Cell(10, 2).Range.Text = "93-93="
$tbl.Cell(10, 3).Range.Text = "81+7="
$tbl.Cell(10, 4).Range.Text = "99-4="
$tbl.Cell(10, 5).Range.Text = "76-75="
$tbl.Cell(11, 1).Range.Text = "19+63="
$tbl.Cell(11, 2).Range.Text = "8+59="
$tbl.Cell(11, 3).Range.Text = "41-4="
$tbl.Cell(11, 4).Range.Text = "84-42="
$tbl.Cell(11, 5).Range.Text = "16+60="
$tbl.Cell(12, 1).Range.Text = "7+63="
$tbl.Cell(12, 2).Range.Text = "21+11="
$tbl.Cell(12, 3).Range.Text = "53-23="
$tbl.Cell(12, 4).Range.Text = "40-21="
$tbl.Cell(12, 5).Range.Text = "27+42="
$tbl.Cell(13, 1).Range.Text = "92-32="
$tbl.Cell(13, 2).Range.Text = "36+34="
$tbl.Cell(13, 3).Range.Text = "94-52="
$tbl.Cell(13, 4).Range.Text = "89-66="
$tbl.Cell(13, 5).Range.Text = "85-10="
$tbl.Cell(14, 1).Range.Text = "18+54="
$tbl.Cell(14, 2).Range.Text = "97-52="
$tbl.Cell(14, 3).Range.Text = "61-59="
$tbl.Cell(14, 4).Range.Text = "50-3="
$tbl.Cell(14, 5).Range.Text = "85+11="
$tbl.Cell(15, 1).Range.Text = "89-6="
$tbl.Cell(15, 2).Range.Text = "7+85="
$tbl.Cell(15, 3).Range.Text = "63+1="
$tbl.Cell(15, 4).Range.Text = "78-54="
$tbl.Cell(15, 5).Range.Text = "13+1="
$tbl.Cell(16, 1).Range.Text = "30-13="
$tbl.Cell(16, 2).Range.Text = "39-31="
$tbl.Cell(16, 3).Range.Text = "46-38="
$tbl.Cell(16, 4).Range.Text = "38+25="
$tbl.Cell(16, 5).Range.Text = "33+64="
$tbl.Cell(17, 1).Range.Text = "24-24="
$tbl.Cell(17, 2).Range.Text = "89-32="
$tbl.Cell(17, 3).Range.Text = "26+61="
$tbl.Cell(17, 4).Range.Text = "45+27="
$tbl.Cell(17, 5).Range.Text = "31+1="
$tbl.Cell(18, 1).Range.Text = "44-17="
$tbl.Cell(18, 2).Range.Text = "47+20="
$tbl.Cell(18, 3).Range.Text = "86-62="
$tbl.Cell(18, 4).Range.Text = "75-60="
$tbl.Cell(18, 5).Range.Text = "46+7="
$tbl.Cell(19, 1).Range.Text = "74-16="
$tbl.Cell(19, 2).Range.Text = "78-36="
$tbl.Cell(19, 3).Range.Text = "11+78="
$tbl.Cell(19, 4).Range.Text = "54-39="
$tbl.Cell(19, 5).Range.Text = "90-61="
$tbl.Cell(20, 1).Range.Text = "86-9="
$tbl.Cell(20, 2).Range.Text = "88-40="
$tbl.Cell(20, 3).Range.Text = "74+1="
$tbl.Cell(20, 4).Range.Text = "8+17="
$tbl.Cell(20, 5).Range.Text = "53+36="
